$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the hours logged for 2017-08-15 (D16): 7 -> 8
$ws.Range("D16").Value = 8

# Add hours logged for 2017-08-16 (D17): new value 8
$ws.Range("D17").Value = 8

# Recalculate so cached formula results (e.g. F3 SUM) update
$excel.Calculate()

# Update the active selection to D17, matching the saved view state
$ws.Range("D17").Select()
